$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously-used range completely (A1:H6) so stale cells (e.g. column H,
# row 6, the old E3 image url) don't linger.
$ws.Range("A1:H6").Clear()

# Force text format on the phone-number and score data cells so the leading
# "+" and exact decimal text are preserved instead of being coerced to
# numbers by Excel's smart-entry parsing.
$ws.Range("A2:A5").NumberFormat = "@"
$ws.Range("D2:D5").NumberFormat = "@"

# Header row
$ws.Range("A1").Value = "Phone Number"
$ws.Range("B1").Value = "ID"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Score"
$ws.Range("E1").Value = "image"
$ws.Range("F1").Value = "Access"
$ws.Range("G1").Value = "city"

# Row 2
$ws.Range("A2").Value = "+919446933709"
$ws.Range("B2").Value = "asmoggEczlaJI9aJKUw3qg=="
$ws.Range("C2").Value = "Cdr Manu Vidyarthi"
$ws.Range("D2").Value = "0.3222276"
$ws.Range("F2").Value = "PUBLIC"
$ws.Range("G2").Value = "Kerala"

# Row 3
$ws.Range("A3").Value = "+919447433709"
$ws.Range("B3").Value = "uaef4zP1ky/K+u5/LAMhUw=="
$ws.Range("C3").Value = "Subani"
$ws.Range("D3").Value = "0.3168635"
$ws.Range("F3").Value = "PUBLIC"
$ws.Range("G3").Value = "Kerala"

# Row 4
$ws.Range("A4").Value = "+918891911108"
$ws.Range("B4").Value = "xdVTVN9NhWSSadKauFRPMg=="
$ws.Range("C4").Value = "Subani Vidyarthi"
$ws.Range("D4").Value = "0.30928558"
$ws.Range("F4").Value = "PUBLIC"
$ws.Range("G4").Value = "Kerala"

# Row 5
$ws.Range("A5").Value = "+919854291183"
$ws.Range("B5").Value = "toIIhLeJJ1wo5K1S6qA7Mw=="
$ws.Range("C5").Value = "Barnali Bhagabati"
$ws.Range("D5").Value = "0.30918905"
$ws.Range("F5").Value = "PUBLIC"
$ws.Range("G5").Value = "Assam"

# Reset selection to A6:A8 to match the committed state.
$ws.Range("A6:A8").Select() | Out-Null
